$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header/shared-string text renames
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# Data value updates (column C: GDP figures recomputed; column AL: Colony flag)
$ws.Range("C2").Value = 2870.311589353206
$ws.Range("C3").Value = 697.6889104500298
$ws.Range("AL3").Value = 1
$ws.Range("C4").Value = 1460.056109840828
$ws.Range("C5").Value = 1909.084588129339
$ws.Range("C6").Value = 10594.98659239237
$ws.Range("C7").Value = 4547.50930098406
$ws.Range("C8").Value = 4132.902312418774
$ws.Range("C9").Value = 1268.249210347625
$ws.Range("C10").Value = 567.9059336271471
$ws.Range("C12").Value = 1299.344949460393
$ws.Range("C13").Value = 1280.225469721551
$ws.Range("C14").Value = 341.5541149051794
$ws.Range("C15").Value = 612.3436990512633
$ws.Range("C16").Value = 2898.942214704482
$ws.Range("C17").Value = 665.6274194933962
$ws.Range("AL17").Value = 1
$ws.Range("C18").Value = 1503.870423231357
$ws.Range("C19").Value = 10385.96443195552
$ws.Range("C20").Value = 1955.461557360978
$ws.Range("C21").Value = 11286.24301624575
$ws.Range("C22").Value = 4633.590358399045
$ws.Range("C23").Value = 4550.453595838572
$ws.Range("C24").Value = 1357.563719132622
$ws.Range("C25").Value = 592.4010974509293
$ws.Range("C27").Value = 1446.371630707023
$ws.Range("C28").Value = 369.2024078290272
$ws.Range("C29").Value = 644.763840173281
$ws.Range("C30").Value = 3083.80337578809
$ws.Range("C31").Value = 2965.153206179127
$ws.Range("C32").Value = 691.8942672110555
$ws.Range("AL32").Value = 1
$ws.Range("C33").Value = 1577.487171555845
$ws.Range("C34").Value = 2024.117324382548
$ws.Range("C35").Value = 11627.81065059172
$ws.Range("C36").Value = 4921.848409120176
$ws.Range("C37").Value = 4961.234688573883
$ws.Range("C38").Value = 1410.426304742003
$ws.Range("C39").Value = 612.032557723897
$ws.Range("C40").Value = 513.7390871590731
$ws.Range("C42").Value = 701.4459636783288
$ws.Range("AL42").Value = 1
$ws.Range("C43").Value = 1657.651524528445
$ws.Range("C44").Value = 2094.024217383061
$ws.Range("C45").Value = 5122.180090208862
$ws.Range("C46").Value = 11745.7759262897
$ws.Range("C47").Value = 642.5204633514562
$ws.Range("C48").Value = 1443.492614888721
$ws.Range("C49").Value = 534.5063430177229
$ws.Range("C51").Value = 720.7128711178943
$ws.Range("AL51").Value = 1
$ws.Range("C52").Value = 1716.389195271215
$ws.Range("C53").Value = 2201.396847776877
$ws.Range("C54").Value = 5295.682695961288
$ws.Range("C55").Value = 11993.48398487312
$ws.Range("C56").Value = 670.8378265228503
$ws.Range("C57").Value = 1505.810948829135
$ws.Range("C58").Value = 951.3148210424945
$ws.Range("C59").Value = 2286.013198234259
$ws.Range("C60").Value = 5412.131646018807
$ws.Range("C61").Value = 449.4203771491282
$ws.Range("C62").Value = 2025.814194788851
$ws.Range("C63").Value = 1640.18070024053
$ws.Range("C64").Value = 558.2093442539386
$ws.Range("C65").Value = 711.3043470146426
$ws.Range("C66").Value = 1775.027517189621
$ws.Range("C68").Value = 1579.189101937001
$ws.Range("C69").Value = 1002.388731936373
$ws.Range("C70").Value = 2361.056581219794
$ws.Range("C71").Value = 5330.539154475424
$ws.Range("C72").Value = 482.6390663355013
$ws.Range("C73").Value = 2067.29003376698
$ws.Range("C74").Value = 1751.664428859304
$ws.Range("C75").Value = 579.0880693780265
$ws.Range("C76").Value = 731.9993357350996
$ws.Range("C77").Value = 1836.014008604312
$ws.Range("C79").Value = 1667.171891046301
$ws.Range("C80").Value = 1062.040157863007
$ws.Range("C81").Value = 2425.561644739583
$ws.Range("C82").Value = 5176.058803160127
$ws.Range("C83").Value = 514.0573067519859
$ws.Range("C84").Value = 2111.193164269742
$ws.Range("C85").Value = 1875.732161108182
$ws.Range("C86").Value = 584.2111078769213
$ws.Range("C87").Value = 729.6614300490079
$ws.Range("C88").Value = 1895.214690888655
$ws.Range("C90").Value = 1752.531946133768
